$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 120
$ws.Range("I12").Value = 90
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = -490

$ws.Range("H62").Value = 26400.428
$ws.Range("I62").Value = 5549.9165
$ws.Range("J62").Value = 151503.5
$ws.Range("K62").Value = 5549.9165
$ws.Range("L62").Value = 151503.5
$ws.Range("M62").Value = -4925.9165
$ws.Range("N62").Value = -152751.5

$ws.Range("H65").Value = 26400.428
$ws.Range("I65").Value = 5549.9165
$ws.Range("J65").Value = 151503.5
$ws.Range("K65").Value = 27749.5825
$ws.Range("L65").Value = 757517.5
$ws.Range("M65").Value = -24629.5825
$ws.Range("N65").Value = -763757.5

$ws.Range("H106").Value = 2344.2144
$ws.Range("I106").Value = 2257.2856
$ws.Range("J106").Value = 2431.1428
$ws.Range("K106").Value = 2257.2856
$ws.Range("L106").Value = 2431.1428
$ws.Range("M106").Value = -1626.2856
$ws.Range("N106").Value = -3693.1428

$ws.Range("H107").Value = 542.7222
$ws.Range("I107").Value = 551.6429000000001
$ws.Range("J107").Value = 511.5
$ws.Range("K107").Value = 551.6429000000001
$ws.Range("L107").Value = 511.5
$ws.Range("M107").Value = 1368.3571
$ws.Range("N107").Value = -4351.5

$ws.Range("H113").Value = 3547.158
$ws.Range("I113").Value = 2390
$ws.Range("J113").Value = 3960.4285
$ws.Range("K113").Value = 2390
$ws.Range("L113").Value = 3960.4285
$ws.Range("M113").Value = 864
$ws.Range("N113").Value = -10468.4285

$ws.Range("H132").Value = 7823.846
$ws.Range("I132").Value = 6917.926
$ws.Range("J132").Value = 9862.166999999999
$ws.Range("K132").Value = 20753.778
$ws.Range("L132").Value = 29586.501
$ws.Range("M132").Value = -18223.778
$ws.Range("N132").Value = -34646.501

$ws.Range("H138").Value = 1827.5
$ws.Range("I138").Value = 1182.75
$ws.Range("J138").Value = 3301.2144
$ws.Range("K138").Value = 3548.25
$ws.Range("L138").Value = 9903.643199999999
$ws.Range("M138").Value = 1591.75
$ws.Range("N138").Value = -20183.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8610.280000000001
$ws.Range("I32").Value = 9018.204
$ws.Range("J32").Value = 6111.75
$ws.Range("K32").Value = 9018.204
$ws.Range("L32").Value = 6111.75
$ws.Range("M32").Value = -8731.204
$ws.Range("N32").Value = -6685.75

$ws.Range("H110").Value = 4311
$ws.Range("I110").Value = 1399.75
$ws.Range("J110").Value = 6640
$ws.Range("K110").Value = 1399.75
$ws.Range("L110").Value = 6640
$ws.Range("M110").Value = 645.25
$ws.Range("N110").Value = -10730

$ws.Range("H132").Value = 5728.6875
$ws.Range("I132").Value = 2562
$ws.Range("J132").Value = 7895.3687
$ws.Range("K132").Value = 7686
$ws.Range("L132").Value = 23686.1061
$ws.Range("M132").Value = -5156
$ws.Range("N132").Value = -28746.1061

$ws.Range("H133").Value = 32980
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 32980
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 32980
$ws.Range("N133").Value = -38040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4452.4165
$ws.Range("I7").Value = 8356.583000000001
$ws.Range("J7").Value = 548.25
$ws.Range("K7").Value = 8356.583000000001
$ws.Range("L7").Value = 548.25
$ws.Range("M7").Value = -8243.583000000001
$ws.Range("N7").Value = -774.25

$ws.Range("H16").Value = 3362.9285
$ws.Range("I16").Value = 3740.1428
$ws.Range("J16").Value = 2985.7144
$ws.Range("K16").Value = 3740.1428
$ws.Range("L16").Value = 2985.7144
$ws.Range("M16").Value = -3453.1428
$ws.Range("N16").Value = -3559.7144

$ws.Range("H26").Value = 5000
$ws.Range("I26").Value = 5000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 5000
$ws.Range("L26").ClearContents()
$ws.Range("M26").Value = -4713
$ws.Range("N26").Value = 0

$ws.Range("H86").Value = 66674972
$ws.Range("I86").Value = 125009750
$ws.Range("J86").Value = 6656.857
$ws.Range("K86").Value = 125009750
$ws.Range("L86").Value = 6656.857
$ws.Range("M86").Value = -125008627
$ws.Range("N86").Value = -8902.857

$ws.Range("H89").Value = 66674972
$ws.Range("I89").Value = 125009750
$ws.Range("J89").Value = 6656.857
$ws.Range("K89").Value = 625048750
$ws.Range("L89").Value = 33284.285
$ws.Range("M89").Value = -625043134
$ws.Range("N89").Value = -44516.285

$ws.Range("H105").Value = 1424.875
$ws.Range("I105").Value = 1442.7142
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1442.7142
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 304.2858000000001
$ws.Range("N105").Value = -4794

$ws.Range("H107").Value = 493.22223
$ws.Range("I107").Value = 348.42856
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 348.42856
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1571.57144
$ws.Range("N107").Value = -4840

$ws.Range("H113").Value = 3362.9285
$ws.Range("I113").Value = 3740.1428
$ws.Range("J113").Value = 2985.7144
$ws.Range("K113").Value = 3740.1428
$ws.Range("L113").Value = 2985.7144
$ws.Range("M113").Value = -1570.1428
$ws.Range("N113").Value = -7325.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 41668100
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 47620640
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 142861920
$ws.Range("M25").Value = -731
$ws.Range("N25").Value = -142862258

$ws.Range("H30").Value = 41668100
$ws.Range("I30").Value = 300
$ws.Range("J30").Value = 47620640
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 142861920
$ws.Range("M30").Value = -798
$ws.Range("N30").Value = -142862124

$ws.Range("H122").Value = 2868.2104
$ws.Range("I122").Value = 498.5
$ws.Range("J122").Value = 3147
$ws.Range("K122").Value = 4486.5
$ws.Range("L122").Value = 28323
$ws.Range("M122").Value = -2036.5
$ws.Range("N122").Value = -33223

$ws.Range("H132").Value = 3655.6897
$ws.Range("I132").Value = 1169.375
$ws.Range("J132").Value = 6715.769
$ws.Range("K132").Value = 10524.375
$ws.Range("L132").Value = 60441.921
$ws.Range("M132").Value = -7994.375
$ws.Range("N132").Value = -65501.921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2659
$ws.Range("I132").Value = 2303.4546
$ws.Range("J132").Value = 2889.0588
$ws.Range("K132").Value = 6910.3638
$ws.Range("L132").Value = 8667.1764
$ws.Range("M132").Value = -4380.3638
$ws.Range("N132").Value = -13727.1764
